$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- D381:D386 bsecode cells: convert stored type from text to numeric ---
$ws.Cells.Item(381, 4).Value = 532827
$ws.Cells.Item(382, 4).Value = 500302
$ws.Cells.Item(383, 4).Value = 500260
$ws.Cells.Item(384, 4).Value = 500085
$ws.Cells.Item(385, 4).Value = 500400
$ws.Cells.Item(386, 4).Value = 532555

# --- Append new rows 387:418 (refreshed "day" timeframe snapshot) ---
$ws.Cells.Item(387, 1).Value = 1
$ws.Cells.Item(387, 2).Value = 'NIFTY'
$ws.Cells.Item(387, 3).Value = 'NIFTY'
$ws.Cells.Item(387, 4).Value = "'"
$ws.Cells.Item(387, 4).Style = "Normal"
$ws.Cells.Item(387, 5).Value = 0.51
$ws.Cells.Item(387, 6).Value = 24698.85
$ws.Cells.Item(387, 7).Value = 0
$ws.Cells.Item(387, 8).Value = 'day'
$ws.Cells.Item(387, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(388, 1).Value = 2
$ws.Cells.Item(388, 2).Value = 'MARUTI'
$ws.Cells.Item(388, 3).Value = 'Maruti Suzuki India Limited'
$ws.Cells.Item(388, 4).Value = "'532500"
$ws.Cells.Item(388, 4).Style = "Normal"
$ws.Cells.Item(388, 5).Value = 0.54
$ws.Cells.Item(388, 6).Value = 12214.95
$ws.Cells.Item(388, 7).Value = 471933
$ws.Cells.Item(388, 8).Value = 'day'
$ws.Cells.Item(388, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(389, 1).Value = 3
$ws.Cells.Item(389, 2).Value = 'OFSS'
$ws.Cells.Item(389, 3).Value = 'Oracle Financial Services Software Limited'
$ws.Cells.Item(389, 4).Value = "'532466"
$ws.Cells.Item(389, 4).Style = "Normal"
$ws.Cells.Item(389, 5).Value = 1.03
$ws.Cells.Item(389, 6).Value = 11056.7
$ws.Cells.Item(389, 7).Value = 222689
$ws.Cells.Item(389, 8).Value = 'day'
$ws.Cells.Item(389, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(390, 1).Value = 4
$ws.Cells.Item(390, 2).Value = 'COFORGE'
$ws.Cells.Item(390, 3).Value = 'Coforge (Niit Tech)'
$ws.Cells.Item(390, 4).Value = "'532541"
$ws.Cells.Item(390, 4).Style = "Normal"
$ws.Cells.Item(390, 5).Value = 0.92
$ws.Cells.Item(390, 6).Value = 6115.25
$ws.Cells.Item(390, 7).Value = 325670
$ws.Cells.Item(390, 8).Value = 'day'
$ws.Cells.Item(390, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(391, 1).Value = 5
$ws.Cells.Item(391, 2).Value = 'LTIM'
$ws.Cells.Item(391, 3).Value = 'LTI Mindtree Ltd'
$ws.Cells.Item(391, 4).Value = "'540005"
$ws.Cells.Item(391, 4).Style = "Normal"
$ws.Cells.Item(391, 5).Value = 0.5600000000000001
$ws.Cells.Item(391, 6).Value = 5707.8
$ws.Cells.Item(391, 7).Value = 326485
$ws.Cells.Item(391, 8).Value = 'day'
$ws.Cells.Item(391, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(392, 1).Value = 6
$ws.Cells.Item(392, 2).Value = 'LTTS'
$ws.Cells.Item(392, 3).Value = 'L&t Technology Services Limited'
$ws.Cells.Item(392, 4).Value = "'540115"
$ws.Cells.Item(392, 4).Style = "Normal"
$ws.Cells.Item(392, 5).Value = 1.74
$ws.Cells.Item(392, 6).Value = 5376.2
$ws.Cells.Item(392, 7).Value = 159878
$ws.Cells.Item(392, 8).Value = 'day'
$ws.Cells.Item(392, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(393, 1).Value = 7
$ws.Cells.Item(393, 2).Value = 'PERSISTENT'
$ws.Cells.Item(393, 3).Value = 'Persistent Systems Limited'
$ws.Cells.Item(393, 4).Value = "'533179"
$ws.Cells.Item(393, 4).Style = "Normal"
$ws.Cells.Item(393, 5).Value = 1.24
$ws.Cells.Item(393, 6).Value = 4935.15
$ws.Cells.Item(393, 7).Value = 390183
$ws.Cells.Item(393, 8).Value = 'day'
$ws.Cells.Item(393, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(394, 1).Value = 8
$ws.Cells.Item(394, 2).Value = 'HAL'
$ws.Cells.Item(394, 3).Value = 'Hindustan Aeronautics Ltd'
$ws.Cells.Item(394, 4).Value = "'541154"
$ws.Cells.Item(394, 4).Style = "Normal"
$ws.Cells.Item(394, 5).Value = -1.17
$ws.Cells.Item(394, 6).Value = 4736.05
$ws.Cells.Item(394, 7).Value = 1472866
$ws.Cells.Item(394, 8).Value = 'day'
$ws.Cells.Item(394, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(395, 1).Value = 9
$ws.Cells.Item(395, 2).Value = 'TCS'
$ws.Cells.Item(395, 3).Value = 'Tata Consultancy Services Limited'
$ws.Cells.Item(395, 4).Value = "'532540"
$ws.Cells.Item(395, 4).Style = "Normal"
$ws.Cells.Item(395, 5).Value = 0.74
$ws.Cells.Item(395, 6).Value = 4523.3
$ws.Cells.Item(395, 7).Value = 2212298
$ws.Cells.Item(395, 8).Value = 'day'
$ws.Cells.Item(395, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(396, 1).Value = 10
$ws.Cells.Item(396, 2).Value = 'INDIGO'
$ws.Cells.Item(396, 3).Value = 'Interglobe Aviation Limited'
$ws.Cells.Item(396, 4).Value = "'539448"
$ws.Cells.Item(396, 4).Style = "Normal"
$ws.Cells.Item(396, 5).Value = 1.66
$ws.Cells.Item(396, 6).Value = 4302.05
$ws.Cells.Item(396, 7).Value = 798389
$ws.Cells.Item(396, 8).Value = 'day'
$ws.Cells.Item(396, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(397, 1).Value = 11
$ws.Cells.Item(397, 2).Value = 'CUMMINSIND'
$ws.Cells.Item(397, 3).Value = 'Cummins India Limited'
$ws.Cells.Item(397, 4).Value = "'500480"
$ws.Cells.Item(397, 4).Style = "Normal"
$ws.Cells.Item(397, 5).Value = 1.99
$ws.Cells.Item(397, 6).Value = 3829.55
$ws.Cells.Item(397, 7).Value = 545141
$ws.Cells.Item(397, 8).Value = 'day'
$ws.Cells.Item(397, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(398, 1).Value = 12
$ws.Cells.Item(398, 2).Value = 'ASIANPAINT'
$ws.Cells.Item(398, 3).Value = 'Asian Paints Limited'
$ws.Cells.Item(398, 4).Value = "'500820"
$ws.Cells.Item(398, 4).Style = "Normal"
$ws.Cells.Item(398, 5).Value = 0.87
$ws.Cells.Item(398, 6).Value = 3103.2
$ws.Cells.Item(398, 7).Value = 668282
$ws.Cells.Item(398, 8).Value = 'day'
$ws.Cells.Item(398, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(399, 1).Value = 13
$ws.Cells.Item(399, 2).Value = 'MPHASIS'
$ws.Cells.Item(399, 3).Value = 'Mphasis Limited'
$ws.Cells.Item(399, 4).Value = "'526299"
$ws.Cells.Item(399, 4).Style = "Normal"
$ws.Cells.Item(399, 5).Value = 1.05
$ws.Cells.Item(399, 6).Value = 3008.95
$ws.Cells.Item(399, 7).Value = 619736
$ws.Cells.Item(399, 8).Value = 'day'
$ws.Cells.Item(399, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(400, 1).Value = 14
$ws.Cells.Item(400, 2).Value = 'RELIANCE'
$ws.Cells.Item(400, 3).Value = 'Reliance Industries Limited'
$ws.Cells.Item(400, 4).Value = "'500325"
$ws.Cells.Item(400, 4).Style = "Normal"
$ws.Cells.Item(400, 5).Value = 0.51
$ws.Cells.Item(400, 6).Value = 2991.9
$ws.Cells.Item(400, 7).Value = 4205904
$ws.Cells.Item(400, 8).Value = 'day'
$ws.Cells.Item(400, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(401, 1).Value = 15
$ws.Cells.Item(401, 2).Value = 'GODREJPROP'
$ws.Cells.Item(401, 3).Value = 'Godrej Properties Limited'
$ws.Cells.Item(401, 4).Value = "'533150"
$ws.Cells.Item(401, 4).Style = "Normal"
$ws.Cells.Item(401, 5).Value = 1.8
$ws.Cells.Item(401, 6).Value = 2978.65
$ws.Cells.Item(401, 7).Value = 715719
$ws.Cells.Item(401, 8).Value = 'day'
$ws.Cells.Item(401, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(402, 1).Value = 16
$ws.Cells.Item(402, 2).Value = 'INDIAMART'
$ws.Cells.Item(402, 3).Value = 'Indiamart Intermesh Ltd'
$ws.Cells.Item(402, 4).Value = "'542726"
$ws.Cells.Item(402, 4).Style = "Normal"
$ws.Cells.Item(402, 5).Value = 1.87
$ws.Cells.Item(402, 6).Value = 2905.75
$ws.Cells.Item(402, 7).Value = 378569
$ws.Cells.Item(402, 8).Value = 'day'
$ws.Cells.Item(402, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(403, 1).Value = 17
$ws.Cells.Item(403, 2).Value = 'HINDUNILVR'
$ws.Cells.Item(403, 3).Value = 'Hindustan Unilever Limited'
$ws.Cells.Item(403, 4).Value = "'500696"
$ws.Cells.Item(403, 4).Style = "Normal"
$ws.Cells.Item(403, 5).Value = 0.31
$ws.Cells.Item(403, 6).Value = 2751.05
$ws.Cells.Item(403, 7).Value = 1097889
$ws.Cells.Item(403, 8).Value = 'day'
$ws.Cells.Item(403, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(404, 1).Value = 18
$ws.Cells.Item(404, 2).Value = 'INFY'
$ws.Cells.Item(404, 3).Value = 'Infosys Limited'
$ws.Cells.Item(404, 4).Value = "'500209"
$ws.Cells.Item(404, 4).Style = "Normal"
$ws.Cells.Item(404, 5).Value = 0.4
$ws.Cells.Item(404, 6).Value = 1872.2
$ws.Cells.Item(404, 7).Value = 3870333
$ws.Cells.Item(404, 8).Value = 'day'
$ws.Cells.Item(404, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(405, 1).Value = 19
$ws.Cells.Item(405, 2).Value = 'HCLTECH'
$ws.Cells.Item(405, 3).Value = 'Hcl Technologies Limited'
$ws.Cells.Item(405, 4).Value = "'532281"
$ws.Cells.Item(405, 4).Style = "Normal"
$ws.Cells.Item(405, 5).Value = 0.49
$ws.Cells.Item(405, 6).Value = 1686.75
$ws.Cells.Item(405, 7).Value = 2786133
$ws.Cells.Item(405, 8).Value = 'day'
$ws.Cells.Item(405, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(406, 1).Value = 20
$ws.Cells.Item(406, 2).Value = 'TECHM'
$ws.Cells.Item(406, 3).Value = 'Tech Mahindra Limited'
$ws.Cells.Item(406, 4).Value = "'532755"
$ws.Cells.Item(406, 4).Style = "Normal"
$ws.Cells.Item(406, 5).Value = 2.13
$ws.Cells.Item(406, 6).Value = 1628.6
$ws.Cells.Item(406, 7).Value = 2442281
$ws.Cells.Item(406, 8).Value = 'day'
$ws.Cells.Item(406, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(407, 1).Value = 21
$ws.Cells.Item(407, 2).Value = 'TATACONSUM'
$ws.Cells.Item(407, 3).Value = 'TATA Consumer Products Ltd'
$ws.Cells.Item(407, 4).Value = "'500800"
$ws.Cells.Item(407, 4).Style = "Normal"
$ws.Cells.Item(407, 5).Value = -0.53
$ws.Cells.Item(407, 6).Value = 1171.2
$ws.Cells.Item(407, 7).Value = 673817
$ws.Cells.Item(407, 8).Value = 'day'
$ws.Cells.Item(407, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(408, 1).Value = 22
$ws.Cells.Item(408, 2).Value = 'JSWSTEEL'
$ws.Cells.Item(408, 3).Value = 'Jsw Steel Limited'
$ws.Cells.Item(408, 4).Value = "'500228"
$ws.Cells.Item(408, 4).Style = "Normal"
$ws.Cells.Item(408, 5).Value = -0.07000000000000001
$ws.Cells.Item(408, 6).Value = 917.15
$ws.Cells.Item(408, 7).Value = 697065
$ws.Cells.Item(408, 8).Value = 'day'
$ws.Cells.Item(408, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(409, 1).Value = 23
$ws.Cells.Item(409, 2).Value = 'BSOFT'
$ws.Cells.Item(409, 3).Value = 'Birlasoft Ltd'
$ws.Cells.Item(409, 4).Value = "'532400"
$ws.Cells.Item(409, 4).Style = "Normal"
$ws.Cells.Item(409, 5).Value = -1.61
$ws.Cells.Item(409, 6).Value = 613.6
$ws.Cells.Item(409, 7).Value = 4544686
$ws.Cells.Item(409, 8).Value = 'day'
$ws.Cells.Item(409, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(410, 1).Value = 24
$ws.Cells.Item(410, 2).Value = 'UPL'
$ws.Cells.Item(410, 3).Value = 'Upl Limited'
$ws.Cells.Item(410, 4).Value = "'512070"
$ws.Cells.Item(410, 4).Style = "Normal"
$ws.Cells.Item(410, 5).Value = 0.99
$ws.Cells.Item(410, 6).Value = 566.15
$ws.Cells.Item(410, 7).Value = 1194020
$ws.Cells.Item(410, 8).Value = 'day'
$ws.Cells.Item(410, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(411, 1).Value = 25
$ws.Cells.Item(411, 2).Value = 'IGL'
$ws.Cells.Item(411, 3).Value = 'Indraprastha Gas Limited'
$ws.Cells.Item(411, 4).Value = "'532514"
$ws.Cells.Item(411, 4).Style = "Normal"
$ws.Cells.Item(411, 5).Value = 0.05
$ws.Cells.Item(411, 6).Value = 548.45
$ws.Cells.Item(411, 7).Value = 846656
$ws.Cells.Item(411, 8).Value = 'day'
$ws.Cells.Item(411, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(412, 1).Value = 26
$ws.Cells.Item(412, 2).Value = 'WIPRO'
$ws.Cells.Item(412, 3).Value = 'Wipro Limited'
$ws.Cells.Item(412, 4).Value = "'507685"
$ws.Cells.Item(412, 4).Style = "Normal"
$ws.Cells.Item(412, 5).Value = 0.9399999999999999
$ws.Cells.Item(412, 6).Value = 524.65
$ws.Cells.Item(412, 7).Value = 4772112
$ws.Cells.Item(412, 8).Value = 'day'
$ws.Cells.Item(412, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(413, 1).Value = 27
$ws.Cells.Item(413, 2).Value = 'VEDL'
$ws.Cells.Item(413, 3).Value = 'Vedanta Limited'
$ws.Cells.Item(413, 4).Value = "'500295"
$ws.Cells.Item(413, 4).Style = "Normal"
$ws.Cells.Item(413, 5).Value = 0.88
$ws.Cells.Item(413, 6).Value = 446.65
$ws.Cells.Item(413, 7).Value = 8887888
$ws.Cells.Item(413, 8).Value = 'day'
$ws.Cells.Item(413, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(414, 1).Value = 28
$ws.Cells.Item(414, 2).Value = 'BANKBARODA'
$ws.Cells.Item(414, 3).Value = 'Bank Of Baroda'
$ws.Cells.Item(414, 4).Value = "'532134"
$ws.Cells.Item(414, 4).Style = "Normal"
$ws.Cells.Item(414, 5).Value = 2.83
$ws.Cells.Item(414, 6).Value = 254.35
$ws.Cells.Item(414, 7).Value = 25268962
$ws.Cells.Item(414, 8).Value = 'day'
$ws.Cells.Item(414, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(415, 1).Value = 29
$ws.Cells.Item(415, 2).Value = 'ABCAPITAL'
$ws.Cells.Item(415, 3).Value = 'Aditya Birla Capital Ltd'
$ws.Cells.Item(415, 4).Value = "'540691"
$ws.Cells.Item(415, 4).Style = "Normal"
$ws.Cells.Item(415, 5).Value = 2.4
$ws.Cells.Item(415, 6).Value = 218.11
$ws.Cells.Item(415, 7).Value = 5285861
$ws.Cells.Item(415, 8).Value = 'day'
$ws.Cells.Item(415, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(416, 1).Value = 30
$ws.Cells.Item(416, 2).Value = 'MOTHERSON'
$ws.Cells.Item(416, 3).Value = 'Motherson Sumi Systems Limited'
$ws.Cells.Item(416, 4).Value = "'517334"
$ws.Cells.Item(416, 4).Style = "Normal"
$ws.Cells.Item(416, 5).Value = 1.82
$ws.Cells.Item(416, 6).Value = 189.51
$ws.Cells.Item(416, 7).Value = 28528687
$ws.Cells.Item(416, 8).Value = 'day'
$ws.Cells.Item(416, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(417, 1).Value = 31
$ws.Cells.Item(417, 2).Value = 'CANBK'
$ws.Cells.Item(417, 3).Value = 'Canara Bank'
$ws.Cells.Item(417, 4).Value = "'532483"
$ws.Cells.Item(417, 4).Style = "Normal"
$ws.Cells.Item(417, 5).Value = 1.33
$ws.Cells.Item(417, 6).Value = 111.36
$ws.Cells.Item(417, 7).Value = 24247713
$ws.Cells.Item(417, 8).Value = 'day'
$ws.Cells.Item(417, 9).Value = '20/08/2024 11:35:40'

$ws.Cells.Item(418, 1).Value = 32
$ws.Cells.Item(418, 2).Value = 'IDFC'
$ws.Cells.Item(418, 3).Value = 'Idfc Limited'
$ws.Cells.Item(418, 4).Value = "'532659"
$ws.Cells.Item(418, 4).Style = "Normal"
$ws.Cells.Item(418, 5).Value = 1.89
$ws.Cells.Item(418, 6).Value = 110.14
$ws.Cells.Item(418, 7).Value = 7628498
$ws.Cells.Item(418, 8).Value = 'day'
$ws.Cells.Item(418, 9).Value = '20/08/2024 11:35:40'

